$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the label in A26 from "has_grimoire" to "zaubern" (new method to create grimoire)
$ws.Range("A26").Value = "zaubern"

# Update the current selection to A26, as reflected in the saved view state
$ws.Range("A26").Select()
